# DynamoDB auth working locally with AWS ignore
#
# Two content edits against the task list ("ListParagraph" / numId=5 items):
#
# 1. "Output feedback the React way only (USING STATE VARIABLES)" —
#    highlight the whole paragraph (incl. paragraph mark) yellow, and
#    append a new, separately-highlighted run: " – Update contact pages??"
#
# 2. "Use the standard way of getting HTML elements (don't use the window
#    object to get elements anymore)." — strike the whole paragraph
#    (incl. paragraph mark).

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {

    $text = $p.Range.Text

    if ($text -like "*Output feedback the React way only*") {

        $para = $p.Range

        # Highlight the existing run *and* the paragraph mark (going through
        # the run's Font object applies the highlight to both the run rPr
        # and the owning paragraph's pPr/rPr, matching how Word records
        # "select whole paragraph + highlight").
        $para.Font.HighlightColorIndex = 7   # wdYellow

        # Insert the follow-up comment as its own run right before the
        # paragraph mark (End - 1), after the highlight above is already in
        # place so it does not get silently merged with the first run.
        $insertAt = $para.End - 1
        $point = $d.Range($insertAt, $insertAt)
        $dash = [char]0x2013
        $newText = " " + $dash + " Update contact pages??"
        $point.InsertAfter($newText)

        # Give the freshly-inserted text its own distinguishable range and
        # force a run split by toggling a format flag on, then back off —
        # this keeps it as a discrete <w:r> instead of re-merging into the
        # previous run even though both end up with identical formatting.
        $newRange = $d.Range($insertAt, $insertAt + $newText.Length)
        $newRange.Bold = 1
        $newRange.Font.HighlightColorIndex = 7   # wdYellow
        $newRange.Bold = 0
    }
    elseif ($text -like "*Use the standard way of getting HTML elements*") {
        # Strike through the run text and the paragraph mark alike.
        $p.Range.Font.StrikeThrough = 1
    }
}
